# Updates the cryptocurrency price/volume table on Sheet1 to refresh
# the scraped market data (prices in column D, 1h volume % in column E),
# matching the "Updated cryptos list ... with GitHub Actions" commit.
#
# All target cells store plain text (t="inlineStr" in the original OOXML),
# e.g. "0.5030", "28.230.57", "  +0.73%  ". Excel's COM layer will
# silently re-interpret a numeric-looking string assigned to .Value as a
# real number (dropping trailing zeros, switching to scientific notation,
# introducing floating point noise, etc.), so we force each target cell
# to Text format ("@") before writing the literal string. This keeps the
# values byte-for-byte identical to what is in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.263.59'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.872.26'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.50'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5030'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3914'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09583'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.77%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.97'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.483'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.00'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.871.85'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.001'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.407'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001129'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.11'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06626'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.27%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.147'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.312.92'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.33'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.281'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.535'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.089.61'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.23'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.81'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.47'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.069'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1053'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.631'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.624'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06744'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.510'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02389'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2178'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.48%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.47'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6354'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.981'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.177'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.002'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.53'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6035'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.14%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.263'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.80'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.71%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.00%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06843'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.23%  '
